# daily auto push: 2026-01-12 22:32 UTC
# Insert one new data row at row 625 (pushing the existing 2026/12/29.. rows
# down by one, all the way to the former last row 666 -> 667), then fill the
# newly opened row with the missing 2026/01/13 06:00 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 625..666 down to 626..667 (also bumps the sheet's used-range
# dimension from A1:D666 to A1:D667 automatically, like a real Excel
# "Insert Sheet Rows").
$ws.Rows.Item(625).Insert()

# Column A holds plain text dates (e.g. "2026/01/13"), not real Excel date
# serials. Force the cell to Text before assigning so Excel's autocorrect
# doesn't reinterpret the string as a date, then drop the temporary format
# again so the cell ends up with no explicit style, matching its neighbours.
$ws.Range("A625").NumberFormat = "@"
$ws.Range("A625").Value = "2026/01/13"
$ws.Range("A625").ClearFormats()

$ws.Range("B625").Value = "火"
$ws.Range("C625").Value = 6
$ws.Range("D625").Value = 201
